# Updates the Bosnia Herzegovina Premier Liga odds sheet:
#  1. Swaps a handful of mismatched match rows (same kickoff time, teams had
#     been entered against the wrong fixture).
#  2. Recomputes the Asian-handicap profit/loss columns PL_Ahh (Y) and
#     PL_Aha (Z) for every match, from FTHG/FTAG, the handicap line (Ah)
#     and the closing Asian-handicap odds (oddAHH/oddAHA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-AhPL {
    param(
        [double]$FTHG,
        [double]$FTAG,
        [double]$Ah,
        [double]$OddH,
        [double]$OddA
    )

    # Quarter handicap lines (e.g. -0.25, -0.75, ...) are really two equally
    # weighted bets split across the two neighbouring half/whole lines.
    $scaled = [Math]::Round($Ah * 4)
    if ([Math]::Abs(($Ah * 4) - $scaled) -gt 0.0000001) {
        $lines = @($Ah)
    } elseif (($scaled % 2) -eq 0) {
        $lines = @($Ah)
    } else {
        $lines = @((($scaled - 1) / 4.0), (($scaled + 1) / 4.0))
    }

    $homeResults = @()
    foreach ($line in $lines) {
        $margin = $FTHG - $FTAG + $line
        if ($margin -gt 0.0000001) { $homeResults += "win" }
        elseif ($margin -lt -0.0000001) { $homeResults += "lose" }
        else { $homeResults += "push" }
    }

    $plH = 0.0
    foreach ($res in $homeResults) {
        if ($res -eq "win") { $plH += ($OddH - 1) }
        elseif ($res -eq "lose") { $plH += -1 }
    }
    $plH = $plH / $homeResults.Count

    $awayResults = @()
    foreach ($res in $homeResults) {
        if ($res -eq "win") { $awayResults += "lose" }
        elseif ($res -eq "lose") { $awayResults += "win" }
        else { $awayResults += "push" }
    }

    $plA = 0.0
    foreach ($res in $awayResults) {
        if ($res -eq "win") { $plA += ($OddA - 1) }
        elseif ($res -eq "lose") { $plA += -1 }
    }
    $plA = $plA / $awayResults.Count

    return @($plH, $plA)
}

# --- 1. Swap the rows that were recorded against the wrong fixture ---------
# Columns B..AB hold all of the match data; column A is just the sequential
# row index and must stay put.
$swapPairs = @(
    @(49, 50),
    @(76, 77),
    @(87, 88),
    @(122, 123)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range("B$r1`:AB$r1")
    $rng2 = $ws.Range("B$r2`:AB$r2")
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

# --- 2. Recompute PL_Ahh (Y) / PL_Aha (Z) for every match row ---------------
$lastRow = $ws.UsedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $fthg = $ws.Cells.Item($r, 7).Value()
    if ($fthg -eq $null) { continue }
    $ftag = $ws.Cells.Item($r, 8).Value()
    $ah = $ws.Cells.Item($r, 16).Value()
    $oddH = $ws.Cells.Item($r, 17).Value()
    $oddA = $ws.Cells.Item($r, 18).Value()

    $pl = Get-AhPL $fthg $ftag $ah $oddH $oddA

    $ws.Cells.Item($r, 25).Value = $pl[0]
    $ws.Cells.Item($r, 26).Value = $pl[1]
}
